# Insert a new weekly price record as the new row 8 of the "Haba" sheet.
# This pushes the former rows 8..105 down to 9..106 (dimension grows to
# A1:R106) and fills the freshly inserted row 8 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 8, shifting existing data down.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the new data point.
$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "Vega Modelo de Temuco"
$ws.Range("C8").Value = "La Araucanía"
$ws.Range("D8").Value = 45257
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 100112026
$ws.Range("G8").Value = "Haba"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 10000
$ws.Range("N8").Value = '$/saco 25 kilos'
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 400
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
